$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy bordered formatting (style index 2) into new rows 8 and 9 ---
$ws.Range("B7:F7").Copy()
$ws.Range("B8:F9").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Set row heights for newly added rows to match existing table rows ---
$ws.Rows.Item(8).RowHeight = 15.75
$ws.Rows.Item(9).RowHeight = 15.75

# --- Row 3 ---
$ws.Range("B3").Value = 'escola/VerificaMedia.py'
$ws.Range("C3").Value = 'Verificar_media'
$ws.Range("D3").Value = 'Enviando uma string'
$ws.Range("E3").Value = '"ola"'
$ws.Range("F3").Value = 'TypeError("É necessário que seja um número, strings não inclusas")'

# --- Row 4 ---
$ws.Range("B4").Value = 'escola/VerificaMedia.py'
$ws.Range("C4").Value = 'Verificar_media'
$ws.Range("D4").Value = 'Enviando um número abaixo de 0'
$ws.Range("E4").Value = -1
$ws.Range("F4").Value = 'ValueError("Grade must be between 0 and 10")'

# --- Row 5 ---
$ws.Range("B5").Value = 'escola/VerificaMedia.py'
$ws.Range("C5").Value = 'Verificar_media'
$ws.Range("D5").Value = 'Enviando um número acima de 11'
$ws.Range("E5").Value = 11
$ws.Range("F5").Value = 'ValueError("Grade must be between 0 and 10")'

# --- Row 6 ---
$ws.Range("B6").Value = 'escola/calcularMedia.py'
$ws.Range("C6").Value = 'calcular_media'
$ws.Range("D6").Value = 'Enviando uma lista vazia'
$ws.Range("E6").Value = '[]'
$ws.Range("F6").Value = 'ValueError("it is not allowed to send an empty list")'

# --- Row 7 ---
$ws.Range("B7").Value = 'escola/calcularMedia.py'
$ws.Range("C7").Value = 'calcular_media'
$ws.Range("D7").Value = 'Enviando uma string'
$ws.Range("E7").Value = '"hi"'
$ws.Range("F7").Value = 'ValueError("it is not allowed to send an empty list")'

# --- Row 8 ---
$ws.Range("B8").Value = 'escola/calcularMedia.py'
$ws.Range("C8").Value = 'calcular_media'
$ws.Range("D8").Value = 'Enviando um número menor que 0'
$ws.Range("E8").Value = '[1.0, -10.0]'
$ws.Range("F8").Value = 'ValueError("grades can be from 0 to 10")'

# --- Row 9 ---
$ws.Range("B9").Value = 'escola/calcularMedia.py'
$ws.Range("C9").Value = 'calcular_media'
$ws.Range("D9").Value = 'Enviando um número maior que 10'
$ws.Range("E9").Value = '[1.0, 11.0]'
$ws.Range("F9").Value = 'ValueError("grades can be from 0 to 10")'

# --- Apply wrap text to F4 (creates new style index 3 matching style 2 + wrapText) ---
$ws.Range("F4").WrapText = $true

# --- Column C width ---
$ws.Columns.Item(3).ColumnWidth = 16.0

# --- Selection ---
$ws.Range("C20").Select()

